$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - flight #11 (Tuesday, Jan 10) FR3696 to Birmingham
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Tuesday, Jan 10"
$ws.Range("C12").Value = "9:25 AM"
$ws.Range("D12").Value = "FR3696"
$ws.Range("E12").Value = "Birmingham"
$ws.Range("F12").Value = "(BHX)"
$ws.Range("G12").Value = "Ryanair "
$ws.Range("H12").Value = "B738"
$ws.Range("I12").Value = "(EI-DLH)"
$ws.Range("J12").Value = "9:21 AM"
$ws.Range("K12").Borders.LineStyle = -4142
$ws.Range("L12").Value = "0 hours, -4 minutes"
$ws.Range("M12").Borders.LineStyle = -4142

# Row 13 - flight #12 (Tuesday, Jan 10) LO3994 to Warsaw
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Tuesday, Jan 10"
$ws.Range("C13").Value = "3:10 PM"
$ws.Range("D13").Value = "LO3994"
$ws.Range("E13").Value = "Warsaw"
$ws.Range("F13").Value = "(WAW)"
$ws.Range("G13").Value = "LOT "
$ws.Range("H13").Value = "E170"
$ws.Range("I13").Value = "(SP-LDF)"
$ws.Range("J13").Value = "2:59 PM"
$ws.Range("K13").Borders.LineStyle = -4142
$ws.Range("L13").Value = "0 hours, -11 minutes"
$ws.Range("M13").Borders.LineStyle = -4142

# Row 14 - flight #13 (Tuesday, Jan 10) W95176 to London
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Tuesday, Jan 10"
$ws.Range("C14").Value = "8:55 PM"
$ws.Range("D14").Value = "W95176"
$ws.Range("E14").Value = "London"
$ws.Range("F14").Value = "(LTN)"
$ws.Range("G14").Value = "Wizz Air "
$ws.Range("H14").Value = "A321"
$ws.Range("I14").Value = "(G-WUKI)"
$ws.Range("J14").Value = "9:04 PM"
$ws.Range("K14").Borders.LineStyle = -4142
$ws.Range("L14").Value = "0 hours, 9 minutes"
$ws.Range("M14").Borders.LineStyle = -4142
